# Auto-generated: apply updated market-price snapshot values
# (scheduled runner refresh) across all 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H32").Value = 6594.9165
$ws.Range("I32").Value = 2528
$ws.Range("J32").Value = 9499.857
$ws.Range("K32").Value = 2528
$ws.Range("L32").Value = 9499.857
$ws.Range("M32").Value = -2202
$ws.Range("N32").Value = -10151.857
$ws.Range("H112").Value = 3042.7334
$ws.Range("I112").Value = 1486.3334
$ws.Range("K112").Value = 4459.0002
$ws.Range("M112").Value = -3351.0002
$ws.Range("H137").Value = 7176.192
$ws.Range("J137").Value = 22072.055
$ws.Range("L137").Value = 66216.16500000001
$ws.Range("N137").Value = -71316.16500000001
$ws.Range("H138").Value = 3257.7727
$ws.Range("J138").Value = 3082.9443
$ws.Range("L138").Value = 9248.832900000001
$ws.Range("N138").Value = -19528.8329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2894.3635
$ws.Range("I45").Value = 2507.5715
$ws.Range("K45").Value = 2507.5715
$ws.Range("M45").Value = -2130.5715
$ws.Range("H74").Value = 16062.434
$ws.Range("I74").Value = 3007.6667
$ws.Range("K74").Value = 3007.6667
$ws.Range("M74").Value = -2133.6667
$ws.Range("H77").Value = 16062.434
$ws.Range("I77").Value = 3007.6667
$ws.Range("K77").Value = 15038.3335
$ws.Range("M77").Value = -10670.3335
$ws.Range("H122").Value = 3183.3901
$ws.Range("I122").Value = 2442.2593
$ws.Range("K122").Value = 7326.777900000001
$ws.Range("M122").Value = -4876.777900000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 855.36365
$ws.Range("I22").Value = 973
$ws.Range("J22").Value = 541.6667
$ws.Range("K22").Value = 973
$ws.Range("L22").Value = 541.6667
$ws.Range("M22").Value = -800
$ws.Range("N22").Value = -887.6667
$ws.Range("H94").Value = 3726.4285
$ws.Range("I94").Value = 2290.6
$ws.Range("K94").Value = 2290.6
$ws.Range("M94").Value = -1839.6
$ws.Range("H105").Value = 3745.4375
$ws.Range("I105").Value = 3728.4666
$ws.Range("K105").Value = 3728.4666
$ws.Range("M105").Value = -1981.4666
$ws.Range("H132").Value = 80779.5
$ws.Range("J132").Value = 80779.5
$ws.Range("L132").Value = 80779.5
$ws.Range("N132").Value = -90899.5
$ws.Range("H134").Value = 9394.611000000001
$ws.Range("I134").Value = 2683.7
$ws.Range("J134").Value = 17783.25
$ws.Range("K134").Value = 8051.099999999999
$ws.Range("L134").Value = 53349.75
$ws.Range("M134").Value = -5516.099999999999
$ws.Range("N134").Value = -58419.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 621.2
$ws.Range("I3").Value = 1225
$ws.Range("J3").Value = 218.66667
$ws.Range("K3").Value = 1225
$ws.Range("L3").Value = 218.66667
$ws.Range("M3").Value = -1112
$ws.Range("N3").Value = -444.66667
$ws.Range("H31").Value = 16865.322
$ws.Range("I31").Value = 2571.611
$ws.Range("K31").Value = 2571.611
$ws.Range("M31").Value = -2276.611
$ws.Range("H34").Value = 16865.322
$ws.Range("I34").Value = 2571.611
$ws.Range("K34").Value = 2571.611
$ws.Range("M34").Value = -2369.611
$ws.Range("H58").Value = 20967.334
$ws.Range("I58").Value = 7626
$ws.Range("J58").Value = 33095.816
$ws.Range("K58").Value = 7626
$ws.Range("L58").Value = 33095.816
$ws.Range("M58").Value = -7423
$ws.Range("N58").Value = -33501.816
$ws.Range("H105").Value = 18685
$ws.Range("I105").Value = 53755
$ws.Range("J105").Value = 1150
$ws.Range("K105").Value = 53755
$ws.Range("L105").Value = 1150
$ws.Range("M105").Value = -52008
$ws.Range("N105").Value = -4644
$ws.Range("H106").Value = 49998.5
$ws.Range("J106").Value = 49999
$ws.Range("L106").Value = 49999
$ws.Range("N106").Value = -52523
$ws.Range("H134").Value = 27783482
$ws.Range("I134").Value = 1577.5834
$ws.Range("J134").Value = 83347300
$ws.Range("K134").Value = 4732.7502
$ws.Range("L134").Value = 250041900
$ws.Range("M134").Value = -2197.7502
$ws.Range("N134").Value = -250046970
$ws.Range("H135").Value = 103666.336
$ws.Range("J135").Value = 103666.336
$ws.Range("L135").Value = 103666.336
$ws.Range("N135").Value = -113806.336
$ws.Range("H136").Value = 20967.334
$ws.Range("I136").Value = 7626
$ws.Range("J136").Value = 33095.816
$ws.Range("K136").Value = 22878
$ws.Range("L136").Value = 99287.448
$ws.Range("M136").Value = -20328
$ws.Range("N136").Value = -104387.448

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 25004264
$ws.Range("I3").Value = 25004264
$ws.Range("K3").Value = 75012792
$ws.Range("M3").Value = -75012680
$ws.Range("H15").Value = 72.42856999999999
$ws.Range("J15").Value = 84
$ws.Range("L15").Value = 252
$ws.Range("N15").Value = -532
$ws.Range("H41").Value = 2122212.2
$ws.Range("J41").Value = 2546614.5
$ws.Range("L41").Value = 7639843.5
$ws.Range("N41").Value = -7640519.5
$ws.Range("H70").Value = 30000
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 30000
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 39434.918
$ws.Range("J20").Value = 39434.918
$ws.Range("L20").Value = 39434.918
$ws.Range("N20").Value = -39924.918
$ws.Range("H39").Value = 18043.818
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 18043.818
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 18043.818
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -19107.818
$ws.Range("H107").Value = 1366.75
$ws.Range("I107").Value = 1237.5834
$ws.Range("J107").Value = 1560.5
$ws.Range("K107").Value = 1237.5834
$ws.Range("L107").Value = 1560.5
$ws.Range("M107").Value = 682.4166
$ws.Range("N107").Value = -5400.5
$ws.Range("H113").Value = 101619.22
$ws.Range("I113").Value = 151512.67
$ws.Range("K113").Value = 151512.67
$ws.Range("M113").Value = -149342.67
$ws.Range("H116").Value = 72000
$ws.Range("J116").Value = 72000
$ws.Range("L116").Value = 72000
$ws.Range("N116").Value = -81178
$ws.Range("H122").Value = 3264.3
$ws.Range("I122").Value = 3164.1428
$ws.Range("J122").Value = 3498
$ws.Range("K122").Value = 9492.428400000001
$ws.Range("L122").Value = 10494
$ws.Range("M122").Value = -7042.428400000001
$ws.Range("N122").Value = -15394
$ws.Range("H132").Value = 34742.15
$ws.Range("I132").Value = 35844.766
$ws.Range("J132").Value = 26472.5
$ws.Range("K132").Value = 107534.298
$ws.Range("L132").Value = 79417.5
$ws.Range("M132").Value = -105004.298
$ws.Range("N132").Value = -84477.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6884.5356
$ws.Range("I22").Value = 6083.1304
$ws.Range("J22").Value = 10571
$ws.Range("K22").Value = 6083.1304
$ws.Range("L22").Value = 10571
$ws.Range("M22").Value = -5788.1304
$ws.Range("N22").Value = -11161
$ws.Range("H27").Value = 6884.5356
$ws.Range("I27").Value = 6083.1304
$ws.Range("J27").Value = 10571
$ws.Range("K27").Value = 6083.1304
$ws.Range("L27").Value = 10571
$ws.Range("M27").Value = -5976.1304
$ws.Range("N27").Value = -10785
$ws.Range("H40").Value = 12777.714
$ws.Range("I40").Value = 13876.667
$ws.Range("J40").Value = 10799.6
$ws.Range("K40").Value = 13876.667
$ws.Range("L40").Value = 10799.6
$ws.Range("M40").Value = -13740.667
$ws.Range("N40").Value = -11071.6
$ws.Range("H46").Value = 387133.7
$ws.Range("I46").Value = 834707.5
$ws.Range("J46").Value = 3499
$ws.Range("K46").Value = 834707.5
$ws.Range("L46").Value = 3499
$ws.Range("M46").Value = -834519.5
$ws.Range("N46").Value = -3875
$ws.Range("H122").Value = 7319.35
$ws.Range("I122").Value = 6865.8335
$ws.Range("J122").Value = 7999.625
$ws.Range("K122").Value = 20597.5005
$ws.Range("L122").Value = 23998.875
$ws.Range("M122").Value = -18147.5005
$ws.Range("N122").Value = -28898.875
$ws.Range("H136").Value = 7556
$ws.Range("I136").Value = 4614.7036
$ws.Range("J136").Value = 11165.772
$ws.Range("K136").Value = 13844.1108
$ws.Range("L136").Value = 33497.31600000001
$ws.Range("M136").Value = -11294.1108
$ws.Range("N136").Value = -38597.31600000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 30395364
$ws.Range("J15").Value = 30395364
$ws.Range("L15").Value = 30395364
$ws.Range("N15").Value = -30395940
$ws.Range("H95").Value = 28786
$ws.Range("J95").Value = 28786
$ws.Range("L95").Value = 28786
$ws.Range("N95").Value = -34278
$ws.Range("H107").Value = 1277.3
$ws.Range("I107").Value = 1404.25
$ws.Range("K107").Value = 4212.75
$ws.Range("M107").Value = -2292.75
$ws.Range("H122").Value = 4761.6665
$ws.Range("I122").Value = 2448.5715
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 7345.7145
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -4895.7145
$ws.Range("N122").Value = -28900
$ws.Range("H126").Value = 20986.105
$ws.Range("J126").Value = 3604.111
$ws.Range("L126").Value = 10812.333
$ws.Range("N126").Value = -15752.333

